# Update employee absence data rows 2-11 with new values per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2;  A=84127; B="Sra. Gabrielly Cavalcante"; C="Marketing";          D="Viagem de negocios"; E=6; F=45079; G=3031.24 }
    @{ Row=3;  A=42960; B="Maria Sophia Vargas";        C="Recursos Humanos";   D="Viagem de negocios"; E=7; F=45083; G=6858.62 }
    @{ Row=4;  A=46946; B="Ana Beatriz Borges";          C="Recursos Humanos";   D="Doenca";             E=6; F=45082; G=5897.21 }
    @{ Row=5;  A=89195; B="Dra. Mariane Rodrigues";      C="Engenharia";         D="Viagem de negocios"; E=6; F=45095; G=5077.7 }
    @{ Row=6;  A=80471; B="Olivia Peixoto";              C="Marketing";          D="Outros";             E=7; F=45083; G=7264.21 }
    @{ Row=7;  A=6139;  B="Emanuel Novaes";              C="Recursos Humanos";   D="Consulta medica";    E=3; F=45096; G=3499.35 }
    @{ Row=8;  A=49873; B="Luiz Fernando da Cruz";       C="Engenharia";         D="Doenca";             E=3; F=45080; G=2564.81 }
    @{ Row=9;  A=15313; B="Dra. Emanuella da Mota";      C="Atendimento ao Cliente"; D="Doenca";         E=5; F=45096; G=7593.97 }
    @{ Row=10; A=12670; B="Helena Costa";                C="Recursos Humanos";   D="Viagem de negocios"; E=3; F=45104; G=8006.13 }
    @{ Row=11; A=22521; B="Anthony Gabriel Costela";     C="Vendas";             D="Viagem de negocios"; E=7; F=45104; G=6980.56 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

$wb.Save()
